$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Descripción" header (A1) to "Nombre"
$ws.Range("A1").Value = "Nombre"

# Add two new columns: "Descripcion" (O) and "Nombre secundario" (P)
$ws.Range("O1").Value = "Descripcion"
$ws.Range("P1").Value = "Nombre secundario"

$ws.Range("O2").Value = "desc 1"
$ws.Range("O3").Value = "desc 2"

$ws.Range("P2").Value = "nombre sec 11"
$ws.Range("P3").Value = "nombre sec 22"

# Set the new "Nombre secundario" column width to match its content (bestFit)
$ws.Columns.Item(16).ColumnWidth = 17.7

# Update the selection to match the new target cell
$ws.Range("O8").Select() | Out-Null
